$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.130.36"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "2.994.03"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "501.29"
$ws.Range("E5").Value = "  -4.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.51"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.32"
$ws.Range("E9").Value = "  -3.92%  "
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.360"
$ws.Range("E11").Value = "  -2.30%  "
$ws.Range("D12").Value = "3.500.08"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.19"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000161"
$ws.Range("E15").Value = "  -4.72%  "
$ws.Range("D16").Value = "57.192.65"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.09"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "2.989.53"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.89"
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.97"
$ws.Range("E21").Value = "  -5.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.492"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.85"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("E26").Value = "  -4.63%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "0.0₃0900"
$ws.Range("E28").Value = "  -6.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.61"
$ws.Range("E29").Value = "  -4.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.11"
$ws.Range("E30").Value = "  -3.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.77"
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("E32").Value = "  -5.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.30"
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.97"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.59"
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.80"
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("E37").Value = "  -6.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.41"
$ws.Range("E38").Value = "  -6.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0668"
$ws.Range("E39").Value = "  -3.23%  "
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").Value = "3.024.25"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.80"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.77"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.643"
$ws.Range("E44").Value = "  -3.15%  "
$ws.Range("D45").Value = "2.204.71"
$ws.Range("E45").Value = "  -5.36%  "
$ws.Range("E46").Value = "  -5.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.97"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.940"
$ws.Range("E48").Value = "  -8.74%  "
$ws.Range("E49").Value = "  -4.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.45"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.81"
$ws.Range("E51").Value = "  -11.45%  "
